$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps its text formatting so values such as
# "1.00" or "3.52" are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '57.189.51'
$ws.Range("E2").Value = '  +7.17%  '
$ws.Range("D3").Value = '3.240.69'
$ws.Range("E3").Value = '  +2.79%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '395.31'
$ws.Range("E5").Value = '  -0.71%  '
$ws.Range("D6").Value = '107.37'
$ws.Range("E6").Value = '  -0.41%  '
$ws.Range("D7").Value = '3.238.02'
$ws.Range("E7").Value = '  +2.85%  '
$ws.Range("D8").Value = '0.570'
$ws.Range("E8").Value = '  +3.77%  '
$ws.Range("E9").Value = '  -0.04%  '
$ws.Range("E10").Value = '  +1.20%  '
$ws.Range("E11").Value = '  +0.10%  '
$ws.Range("D12").Value = '0.0967'
$ws.Range("E12").Value = '  +11.13%  '
$ws.Range("E13").Value = '  +1.54%  '
$ws.Range("D14").Value = '3.754.27'
$ws.Range("E14").Value = '  +2.81%  '
$ws.Range("D15").Value = '8.11'
$ws.Range("E15").Value = '  +1.63%  '
$ws.Range("D16").Value = '18.91'
$ws.Range("E16").Value = '  -0.64%  '
$ws.Range("D17").Value = '3.259.91'
$ws.Range("E17").Value = '  +3.50%  '
$ws.Range("E18").Value = '  -2.44%  '
$ws.Range("D19").Value = '11.01'
$ws.Range("E19").Value = '  +3.95%  '
$ws.Range("D20").Value = '56.974.62'
$ws.Range("E20").Value = '  +6.81%  '
$ws.Range("E21").Value = '  +0.92%  '
$ws.Range("D22").Value = '0.0000106'
$ws.Range("E22").Value = '  +8.84%  '
$ws.Range("D23").Value = '12.97'
$ws.Range("E23").Value = '  +1.04%  '
$ws.Range("D24").Value = '298.57'
$ws.Range("E24").Value = '  +10.19%  '
$ws.Range("D25").Value = '74.06'
$ws.Range("E25").Value = '  +4.46%  '
$ws.Range("D26").Value = '3.14'
$ws.Range("E26").Value = '  -2.19%  '
$ws.Range("E27").Value = '  +3.14%  '
$ws.Range("D28").Value = '27.88'
$ws.Range("E28").Value = '  +0.49%  '
$ws.Range("D29").Value = '7.68'
$ws.Range("E29").Value = '  -5.19%  '
$ws.Range("E30").Value = '  -2.82%  '
$ws.Range("D31").Value = '0.169'
$ws.Range("E31").Value = '  -0.41%  '
$ws.Range("E32").Value = '  -0.08%  '
$ws.Range("D33").Value = '0.109'
$ws.Range("E33").Value = '  -1.19%  '
$ws.Range("D34").Value = '11.00'
$ws.Range("E34").Value = '  -0.01%  '
$ws.Range("D35").Value = '37.43'
$ws.Range("E35").Value = '  +0.53%  '
$ws.Range("E36").Value = '  -2.36%  '
$ws.Range("E37").Value = '  +1.27%  '
$ws.Range("D38").Value = '51.71'
$ws.Range("E38").Value = '  +2.64%  '
$ws.Range("B39").Value = 'LidoDAOToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D39").Value = '3.52'
$ws.Range("E39").Value = '  -2.59%  '
$ws.Range("B40").Value = 'FirstDigitalUSD'
$ws.Range("C40").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D40").Value = '0.999'
$ws.Range("E40").Value = '  -0.15%  '
$ws.Range("D41").Value = '3.03'
$ws.Range("E41").Value = '  +9.09%  '
$ws.Range("D42").Value = '134.70'
$ws.Range("E42").Value = '  +3.24%  '
$ws.Range("B43").Value = 'Stellar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D43").Value = '0.120'
$ws.Range("E43").Value = '  +2.01%  '
$ws.Range("B44").Value = 'ARBITRUM'
$ws.Range("C44").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D44").Value = '1.90'
$ws.Range("E44").Value = '  -0.42%  '
$ws.Range("B45").Value = 'NEARProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D45").Value = '3.96'
$ws.Range("E45").Value = '  -4.44%  '
$ws.Range("D46").Value = '16.89'
$ws.Range("E46").Value = '  -3.13%  '
$ws.Range("D47").Value = '0.281'
$ws.Range("E47").Value = '  -4.32%  '
$ws.Range("D48").Value = '21.88'
$ws.Range("E48").Value = '  -2.20%  '
$ws.Range("D49").Value = '2.147.58'
$ws.Range("E49").Value = '  +2.55%  '
$ws.Range("D50").Value = '2.09'
$ws.Range("E50").Value = '  +0.83%  '
$ws.Range("E51").Value = '  +24.00%  '

# Restore the default (unstyled) cell style now that the text values are set,
# matching the original workbook formatting.
$ws.Range("D2:D51").Style = "Normal"
